$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").Value = "Graphical data analysis"
$ws.Range("D5").Select()
